# Replace each "m:..." Word field (fldChar begin / instrText / fldChar end)
# with plain literal text runs "{<code>}" -- mirrors the
# TokenIteratorFieldRewriterSplit behaviour: the field delimiter spaces are
# dropped, the code is wrapped in curly braces, and -- for fields whose
# instrText was originally split across several runs -- the same run split
# points are preserved in the resulting literal text.

$d = $word.ActiveDocument

# Known internal run-split points (character offsets into the *trimmed*
# field code, i.e. after stripping the single leading/trailing delimiter
# space) for the two fields whose instrText was split across multiple runs
# in the source document.
$splitMap = @{
    "m:userdoc 'zone2'" = @(1, 16)
    "m:userdoc 'zone3'" = @(15, 16)
}

$bookmarkCounter = 0

while ($d.Fields.Count -gt 0) {
    $f = $d.Fields.Item(1)

    $code = $f.Code.Text
    # Field codes here always carry exactly one leading and one trailing
    # delimiter space (" m:userdoc 'zone1' " etc.) -- strip both.
    $trimmed = $code.Substring(1, $code.Length - 2)
    $insertText = "{" + $trimmed + "}"

    # Position of the field's "begin" character -- this is where the
    # replacement literal text must be inserted.
    $startPos = $f.Code.Start - 1

    $f.Delete()

    # Use InsertBefore (not InsertAfter) on this zero-width range: when a
    # hidden bookmark (e.g. Word's "_GoBack") sits exactly at $startPos,
    # InsertBefore places the new text ahead of it, keeping the bookmark
    # trailing the text -- matching its original position right after the
    # field it used to follow.
    $rng = $d.Range($startPos, $startPos)
    $rng.InsertBefore($insertText)

    # If this field's code used to be split across several runs, recreate
    # that split in the new literal-text run using temporary bookmarks:
    # inserting (then removing) a bookmark at a position forces the run to
    # be split there, and the split survives the bookmark's removal.
    if ($splitMap.ContainsKey($trimmed)) {
        $offsets = $splitMap[$trimmed]
        $names = @()
        foreach ($off in $offsets) {
            $bookmarkCounter = $bookmarkCounter + 1
            $name = "m2docTmpSplit" + $bookmarkCounter
            $names = $names + $name
            $bmPos = $startPos + 1 + $off
            $d.Bookmarks.Add($name, $d.Range($bmPos, $bmPos))
        }
        foreach ($name in $names) {
            $d.Bookmarks.Item($name).Delete()
        }
    }
}
